# The workbook lists distinct "good" labels (column A) together with their
# occurrence counts (column B), sorted descending by count. Re-running the
# notebook that produced this sheet recomputed the same (label, count) pairs,
# but ties in the count column were broken in a different order than before
# (e.g. a different Python dict/set iteration order), so some adjacent rows
# that share an identical count swapped which label sits on which row.
# The counts themselves, and the row positions they occupy, are unchanged -
# only which label is attached to a handful of tied-count rows differs.
#
# Updating the label text for the affected rows causes the engine to
# regenerate the shared-string table to match, which is the only observable
# change in the underlying XML (sharedStrings.xml reordered; everything else
# - dimensions, counts, styles - stays the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "особливый товар"
$ws.Range("A16").Value = "мелочь"

$ws.Range("A19").Value = "небогатый товар"
$ws.Range("A20").Value = "крамными товар"

$ws.Range("A25").Value = "пушной товар"
$ws.Range("A26").Value = "нужный товар"

$ws.Range("A28").Value = "суровский товар"
$ws.Range("A29").Value = "недорогой товар"
$ws.Range("A30").Value = "медный товар"

$ws.Range("A31").Value = "внутренний товар"
$ws.Range("A32").Value = "питейный припасы"

$ws.Range("A35").Value = "заморский товар"
$ws.Range("A36").Value = "купецкий товар"
$ws.Range("A37").Value = "произрастание"
$ws.Range("A38").Value = "галантерейный товар"

$ws.Range("A39").Value = "рукодельный товар"
$ws.Range("A40").Value = "меховой товар"
$ws.Range("A41").Value = "домовый товар"
$ws.Range("A42").Value = "надлежащий товар"
$ws.Range("A43").Value = "харчевой припасы"
